$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.519.57'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '3.445.43'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = '144.77'
$ws.Range('E6').Value = '  +6.55%  '
$ws.Range('D7').Value = '3.446.83'
$ws.Range('E7').Value = '  +2.23%  '
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').Value = '4.037.78'
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '28.13'
$ws.Range('E14').Value = '  +8.07%  '
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').Value = '3.448.14'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = '61.659.09'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = '6.26'
$ws.Range('E19').Value = '  +7.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.20'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.37%  '
$ws.Range('D21').Value = '9.51'
$ws.Range('E21').Value = '  +3.14%  '
$ws.Range('D22').Value = '394.98'
$ws.Range('E22').Value = '  +6.37%  '
$ws.Range('E23').Value = '  +3.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.40%  '
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('D28').Value = '3.591.49'
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('D30').Value = '7.59'
$ws.Range('E30').Value = '  +4.06%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '8.15'
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('E33').Value = '  -8.68%  '
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '23.98'
$ws.Range('E36').Value = '  +3.18%  '
$ws.Range('D37').Value = '3.478.89'
$ws.Range('E37').Value = '  +2.55%  '
$ws.Range('D38').Value = '7.01'
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').Value = '167.59'
$ws.Range('E41').Value = '  +1.77%  '
$ws.Range('D42').Value = '0.0781'
$ws.Range('E42').Value = '  +3.12%  '
$ws.Range('E43').Value = '  +11.42%  '
$ws.Range('D44').Value = '0.803'
$ws.Range('E44').Value = '  +4.29%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '4.49'
$ws.Range('E47').Value = '  +4.36%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '42.25'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').Value = '2.602.64'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').Value = '1.16'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.90'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.44%  '
